# Converts an "RRGGBB" hex string into the BGR-ordered integer that the
# Excel object model's Interior.Color / Font.Color expect.
function RgbColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$xlNone = -4142

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting")

# --- Clear out the old D1:I3 swatch/label block so we can rebuild it ---
$ws.Range("D1:I3").ClearContents()
$ws.Range("D1:I3").Interior.Pattern = $xlNone
$ws.Range("D1:I3").Font.Name = "Calibri"

# --- Row 1: "Custom" palette (new user-picked colors) ---
$ws.Range("D1").Value = "Custom"
$ws.Range("D1").Font.Name = "Calibri (Body)"
$ws.Range("E1").Interior.Color = RgbColor("BB4444")
$ws.Range("F1").Interior.Color = RgbColor("EE9988")
$ws.Range("G1").Interior.Color = RgbColor("FFFFFF")
$ws.Range("H1").Interior.Color = RgbColor("77AADD")
$ws.Range("I1").Interior.Color = RgbColor("4477AA")

# --- Row 2: "Tropic" palette (previously the "Default" swatches) ---
$ws.Range("D2").Value = "Tropic"
$ws.Range("E2").Interior.Color = RgbColor("019C9F")
$ws.Range("F2").Interior.Color = RgbColor("78C6C7")
$ws.Range("G2").Interior.Color = RgbColor("F1F1F1")
$ws.Range("H2").Interior.Color = RgbColor("DFA8CA")
$ws.Range("I2").Interior.Color = RgbColor("C95BA7")

# --- Row 3: "Green-Orange" palette (unchanged colors, now labeled in D) ---
$ws.Range("D3").Value = "Green-Orange"
$ws.Range("E3").Interior.Color = RgbColor("00C74B")
$ws.Range("F3").Interior.Color = RgbColor("94D69F")
$ws.Range("G3").Interior.Color = RgbColor("E2E2E2")
$ws.Range("H3").Interior.Color = RgbColor("F3BB98")
$ws.Range("I3").Interior.Color = RgbColor("F4952B")

# --- Data validation list now reads from the label column D instead of I ---
$dv = $ws.Range("B1").Validation
$dv.Modify($dv.Type, $dv.AlertStyle, $dv.Operator, "=`$D`$1:`$D`$3")

# --- Selected palette becomes "Tropic" ---
$ws.Range("B1").Value = "Tropic"

$ws.Range("B1").Select()
